$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Minoh_room")
$ws.Name = "Minoh_room23"
